$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 29476.25
$ws.Range("J116").Value = 5633.3335
$ws.Range("L116").Value = 5633.3335
$ws.Range("N116").Value = -12517.3335
$ws.Range("H137").Value = 3713743.5
$ws.Range("I137").Value = 7605.6665
$ws.Range("J137").Value = 8346415.5
$ws.Range("K137").Value = 22816.9995
$ws.Range("L137").Value = 25039246.5
$ws.Range("M137").Value = -20266.9995
$ws.Range("N137").Value = -25044346.5
$ws.Range("H141").Value = 3654.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1526843.8
$ws.Range("I61").Value = 51942.953
$ws.Range("J61").Value = 3329500.2
$ws.Range("K61").Value = 51942.953
$ws.Range("L61").Value = 3329500.2
$ws.Range("M61").Value = -51730.953
$ws.Range("N61").Value = -3329924.2
$ws.Range("H63").Value = 16904.387
$ws.Range("I63").Value = 4473
$ws.Range("J63").Value = 23741.65
$ws.Range("K63").Value = 4473
$ws.Range("L63").Value = 23741.65
$ws.Range("M63").Value = -3787
$ws.Range("N63").Value = -25113.65
$ws.Range("H66").Value = 16904.387
$ws.Range("I66").Value = 4473
$ws.Range("J66").Value = 23741.65
$ws.Range("K66").Value = 22365
$ws.Range("L66").Value = 118708.25
$ws.Range("M66").Value = -18933
$ws.Range("N66").Value = -125572.25
$ws.Range("H74").Value = 647970.1
$ws.Range("I74").Value = 2868.8823
$ws.Range("K74").Value = 2868.8823
$ws.Range("M74").Value = -1994.8823
$ws.Range("H77").Value = 647970.1
$ws.Range("I77").Value = 2868.8823
$ws.Range("K77").Value = 14344.4115
$ws.Range("M77").Value = -9976.411500000002
$ws.Range("H110").Value = 4068.923
$ws.Range("I110").Value = 2566.8333
$ws.Range("J110").Value = 5356.4287
$ws.Range("K110").Value = 2566.8333
$ws.Range("L110").Value = 5356.4287
$ws.Range("M110").Value = -521.8332999999998
$ws.Range("N110").Value = -9446.4287
$ws.Range("H122").Value = 608.4091
$ws.Range("I122").Value = 384.93332
$ws.Range("K122").Value = 1154.79996
$ws.Range("M122").Value = 1295.20004
$ws.Range("H136").Value = 1526843.8
$ws.Range("I136").Value = 51942.953
$ws.Range("J136").Value = 3329500.2
$ws.Range("K136").Value = 155828.859
$ws.Range("L136").Value = 9988500.600000001
$ws.Range("M136").Value = -153278.859
$ws.Range("N136").Value = -9993600.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7781.636
$ws.Range("J86").Value = 9631.416999999999
$ws.Range("L86").Value = 9631.416999999999
$ws.Range("N86").Value = -11877.417
$ws.Range("H89").Value = 7781.636
$ws.Range("J89").Value = 9631.416999999999
$ws.Range("L89").Value = 48157.085
$ws.Range("N89").Value = -59389.085
$ws.Range("H94").Value = 977.3125
$ws.Range("I94").Value = 1105
$ws.Range("J94").Value = 287.8
$ws.Range("K94").Value = 1105
$ws.Range("L94").Value = 287.8
$ws.Range("M94").Value = -654
$ws.Range("N94").Value = -1189.8
$ws.Range("H99").Value = 12537.117
$ws.Range("I99").Value = 11598.643
$ws.Range("J99").Value = 16916.666
$ws.Range("K99").Value = 11598.643
$ws.Range("L99").Value = 16916.666
$ws.Range("M99").Value = -10100.643
$ws.Range("N99").Value = -19912.666
$ws.Range("H105").Value = 12425
$ws.Range("I105").Value = 10271.667
$ws.Range("K105").Value = 10271.667
$ws.Range("M105").Value = -8524.666999999999
$ws.Range("H134").Value = 18368794
$ws.Range("I134").Value = 1297.5588
$ws.Range("K134").Value = 3892.6764
$ws.Range("M134").Value = -1357.6764

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3107.0688
$ws.Range("I31").Value = 3785.077
$ws.Range("K31").Value = 3785.077
$ws.Range("M31").Value = -3490.077
$ws.Range("H34").Value = 3107.0688
$ws.Range("I34").Value = 3785.077
$ws.Range("K34").Value = 3785.077
$ws.Range("M34").Value = -3583.077
$ws.Range("H58").Value = 1857.2413
$ws.Range("I58").Value = 1078.2632
$ws.Range("J58").Value = 3337.3
$ws.Range("K58").Value = 1078.2632
$ws.Range("L58").Value = 3337.3
$ws.Range("M58").Value = -875.2632000000001
$ws.Range("N58").Value = -3743.3
$ws.Range("H99").Value = 20000000
$ws.Range("I99").Value = 20000000
$ws.Range("K99").Value = 20000000
$ws.Range("M99").Value = -19998502
$ws.Range("H107").Value = 1758.5333
$ws.Range("J107").Value = 1832.25
$ws.Range("L107").Value = 1832.25
$ws.Range("N107").Value = -5672.25
$ws.Range("H126").Value = 20000000
$ws.Range("I126").Value = 20000000
$ws.Range("K126").Value = 60000000
$ws.Range("M126").Value = -59997530
$ws.Range("H130").Value = 137995.8
$ws.Range("J130").Value = 137995.8
$ws.Range("L130").Value = 137995.8
$ws.Range("N130").Value = -148035.8
$ws.Range("H134").Value = 2998.75
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 2998.75
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").Value = 8996.25
$ws.Range("N134").Value = -14066.25
$ws.Range("H136").Value = 1857.2413
$ws.Range("I136").Value = 1078.2632
$ws.Range("J136").Value = 3337.3
$ws.Range("K136").Value = 3234.7896
$ws.Range("L136").Value = 10011.9
$ws.Range("M136").Value = -684.7896000000001
$ws.Range("N136").Value = -15111.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 33045.176
$ws.Range("I44").Value = 600
$ws.Range("J44").Value = 39997.715
$ws.Range("K44").Value = 1800
$ws.Range("L44").Value = 119993.145
$ws.Range("M44").Value = -1402
$ws.Range("N44").Value = -120789.145
$ws.Range("H131").Value = 5052775.5
$ws.Range("I131").Value = 9092302
$ws.Range("K131").Value = 27276906
$ws.Range("M131").Value = -27271866
$ws.Range("H138").Value = 3197.8333
$ws.Range("I138").Value = 3137.4
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 9412.200000000001
$ws.Range("L138").Value = 10500
$ws.Range("M138").Value = -4272.200000000001
$ws.Range("N138").Value = -20780
$ws.Range("H139").Value = 2708.037
$ws.Range("I139").Value = 1874.1305
$ws.Range("K139").Value = 5622.3915
$ws.Range("M139").Value = -482.3914999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 45753.332
$ws.Range("J96").Value = 45753.332
$ws.Range("L96").Value = 45753.332
$ws.Range("N96").Value = -51245.332
$ws.Range("H100").Value = 37199.8
$ws.Range("I100").Value = 29999
$ws.Range("J100").Value = 39000
$ws.Range("K100").Value = 29999
$ws.Range("L100").Value = 39000
$ws.Range("M100").Value = -28917
$ws.Range("N100").Value = -41164
$ws.Range("H106").Value = 38998.668
$ws.Range("J106").Value = 38998.668
$ws.Range("L106").Value = 38998.668
$ws.Range("N106").Value = -41522.668
$ws.Range("H113").Value = 6900
$ws.Range("I113").Value = 6900
$ws.Range("K113").Value = 6900
$ws.Range("M113").Value = -4730
$ws.Range("H122").Value = 3836.2856
$ws.Range("I122").Value = 2347.4666
$ws.Range("K122").Value = 7042.399800000001
$ws.Range("M122").Value = -4592.399800000001
$ws.Range("H132").Value = 1071456.5
$ws.Range("I132").Value = 2435.3125
$ws.Range("J132").Value = 3209498.8
$ws.Range("K132").Value = 7305.9375
$ws.Range("L132").Value = 9628496.399999999
$ws.Range("M132").Value = -4775.9375
$ws.Range("N132").Value = -9633556.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10227.823
$ws.Range("I7").Value = 6731
$ws.Range("J7").Value = 11303.77
$ws.Range("K7").Value = 6731
$ws.Range("L7").Value = 11303.77
$ws.Range("M7").Value = -6619
$ws.Range("N7").Value = -11527.77
$ws.Range("H88").Value = 10666
$ws.Range("J88").Value = 14999
$ws.Range("L88").Value = 14999
$ws.Range("N88").Value = -15855
$ws.Range("H91").Value = 10666
$ws.Range("J91").Value = 14999
$ws.Range("L91").Value = 14999
$ws.Range("N91").Value = -17963
$ws.Range("H126").Value = 10227.823
$ws.Range("I126").Value = 6731
$ws.Range("J126").Value = 11303.77
$ws.Range("K126").Value = 20193
$ws.Range("L126").Value = 33911.31
$ws.Range("M126").Value = -17723
$ws.Range("N126").Value = -38851.31

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2859786
$ws.Range("I107").Value = 2369.25
$ws.Range("J107").Value = 4764730.5
$ws.Range("K107").Value = 7107.75
$ws.Range("L107").Value = 14294191.5
$ws.Range("M107").Value = -5187.75
$ws.Range("N107").Value = -14298031.5
$ws.Range("H132").Value = 2156.2
$ws.Range("J132").Value = 2859.2727
$ws.Range("L132").Value = 8577.8181
$ws.Range("N132").Value = -13637.8181
$ws.Range("H136").Value = 699.1818
$ws.Range("I136").Value = 556.35297
$ws.Range("K136").Value = 1669.05891
$ws.Range("M136").Value = 880.9410899999998
